$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Valor Mora" values between the first ("2104") and third ("2106")
# periodo-mora rows. Row 17 ("2105") keeps its existing value.
$ws.Range("F16").Value = 39480
$ws.Range("F18").Value = 30268
